$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("F14").Value = 1
$ws.Range("N14").Value = -78.947368421052
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -20
$ws.Range("I16").Value = 62
$ws.Range("J16").Value = 47
$ws.Range("K16").Value = 31.914893617021
$ws.Range("L16").Value = -15.068493150684
$ws.Range("M16").Value = -75.100401606425
$ws.Range("N16").Value = -92.269326683291
$ws.Range("C17").Value = 8
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -18.181818181818
$ws.Range("I17").Value = 180
$ws.Range("J17").Value = 205
$ws.Range("K17").Value = -12.195121951219
$ws.Range("L17").Value = -9.547738693467
$ws.Range("M17").Value = -18.181818181818
$ws.Range("N17").Value = -44.099378881987
$ws.Range("D18").Value = 3
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -80
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 12.222222222222
$ws.Range("L18").Value = -16.528925619834
$ws.Range("M18").Value = -65.292096219931
$ws.Range("N18").Value = -91.990483743061
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 16
$ws.Range("H19").Value = 131.25
$ws.Range("I19").Value = 310
$ws.Range("J19").Value = 267
$ws.Range("K19").Value = 16.104868913857
$ws.Range("L19").Value = -7.462686567164
$ws.Range("M19").Value = -12.429378531073
$ws.Range("N19").Value = -35.281837160751
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = -48
$ws.Range("I20").Value = 147
$ws.Range("J20").Value = 210
$ws.Range("K20").Value = -30
$ws.Range("L20").Value = 12.213740458015
$ws.Range("M20").Value = -51.803278688524
$ws.Range("N20").Value = -94.440242057488
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -5
$ws.Range("F21").Value = 75
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -3.846153846153
$ws.Range("I21").Value = 818
$ws.Range("J21").Value = 834
$ws.Range("K21").Value = -1.918465227817
$ws.Range("L21").Value = -5.977011494252
$ws.Range("M21").Value = -43.469246717346
$ws.Range("N21").Value = -85.287769784172
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 8.333333333333
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 13.793103448275
$ws.Range("I24").Value = 536
$ws.Range("J24").Value = 542
$ws.Range("K24").Value = -1.107011070110
$ws.Range("L24").Value = -17.156105100463
$ws.Range("M24").Value = -19.760479041916
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 130
$ws.Range("J25").Value = 133
$ws.Range("K25").Value = -2.255639097744
$ws.Range("L25").Value = 7.438016528925
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 28
$ws.Range("H26").Value = -3.448275862068
$ws.Range("I26").Value = 312
$ws.Range("J26").Value = 328
$ws.Range("K26").Value = -4.878048780487
$ws.Range("L26").Value = 11.428571428571
$ws.Range("M26").Value = -36.326530612244
$ws.Range("F28").Value = 2
$ws.Range("I28").Value = 26
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = -13.333333333333
$ws.Range("L28").Value = 36.842105263157
$ws.Range("N29").Value = -98.333333333333
$ws.Range("N30").Value = -98.113207547169

# --- Cells changing between numeric and text styles ---
# Row 15 / 27: numeric -> text "0" / "***.*" (style 14/15 -> 13)
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "0"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "0"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "***.*"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100

# --- Fix styles after text/number conversions to match target cellXfs ---
$ws.Range("C33").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("F33").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("H33").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("H28").PasteSpecial(-4122)

# --- Column H width fix (bestFit recalculated by Excel after data change) ---
$ws.Columns.Item(8).ColumnWidth = 6.168446

$excel.CutCopyMode = $false
